{"js": "const body = context.document.body;\n\n// --- Change 1: split the \"Crit\u00e9rio\" run so the M1 equation starts on its own line ---\nconst criterioMatches = body.search(\"M1 = 0,6*NR + 0,4* NP\", { matchCase: true });\ncriterioMatches.load(\"text\");\nawait context.sync();\n\nif (criterioMatches.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the M1 equation text, found \" + criterioMatches.items.length);\n}\n\n// Prefix the formula with a manual line break (appears as U+000B in the Word text model)\n// so the OOXML gains a <w:br/> splitting the run text into two <w:t> runs of text.\ncriterioMatches.items[0].insertText(\"\\u000bM1 = 0,6*NR + 0,4* NP\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Change 2: split the Bibliografia paragraph's single run into four, one per reference ---\nbody.paragraphs.load(\"text\");\nawait context.sync();\n\nconst biblioIndex = body.paragraphs.items.findIndex((p) => p.text.startsWith(\"1)FOX, R.W.\"));\nif (biblioIndex === -1) {\n  throw new Error(\"Could not find the Bibliografia paragraph\");\n}\n\nconst biblioParagraph = body.paragraphs.items[biblioIndex];\nconst references = [\n  \"1)FOX, R.W.; PRITCHARD, P.J.; McDONALD, A.T. Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos Fluidos. Ed. Gen LTC, 7 ed, Rio de Janeiro/RJ, 2010.\",\n  \"2)\u00c7ENGEL, Y.A.; COMBALA, J.M. Mec\u00e2nica dos Fluidos: fundamentos e aplica\u00e7\u00f5es. McGraw-Hill Education (AMGH Editora Ltda),  Porto Alegre/ RS, 2007.\",\n  \"3)COUPER, JR.; PENNEY, W.R.; FAIR, J.R.; WALAS, S.M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005.\",\n  \"4)TROPEA, C.; YARIN, A.L.; FOSS, J.F. Handbook of Experimental Fluid Mechanics. Ed Springer. Springer-Verlag Berlin Heidelberg. 2007\",\n];\nbiblioParagraph.getRange().insertText(references.join(\"\\u000b\"), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: split the \"Crit\u00e9rio\" run so the M1 equation starts on its own line ---\n$find = $d.Content\n$found = $find.Find.Execute(\"M1 = 0,6*NR + 0,4* NP\")\nif (-not $found) {\n    throw \"Could not find the M1 equation text\"\n}\n# Prefixing with a manual line break character (Chr 11, the Word \"vertical tab\"\n# used for <w:br/>) splits the run's text into two <w:t> runs joined by <w:br/>.\n$find.Text = [char]11 + \"M1 = 0,6*NR + 0,4* NP\"\n\n# --- Change 2: split the Bibliografia paragraph's single run into four, one per reference ---\n$biblioIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.StartsWith(\"1)FOX, R.W.\")) {\n        $biblioIndex = $i\n        break\n    }\n}\nif ($biblioIndex -eq -1) {\n    throw \"Could not find the Bibliografia paragraph\"\n}\n\n$references = @(\n    \"1)FOX, R.W.; PRITCHARD, P.J.; McDONALD, A.T. Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos Fluidos. Ed. Gen LTC, 7 ed, Rio de Janeiro/RJ, 2010.\",\n    \"2)\u00c7ENGEL, Y.A.; COMBALA, J.M. Mec\u00e2nica dos Fluidos: fundamentos e aplica\u00e7\u00f5es. McGraw-Hill Education (AMGH Editora Ltda),  Porto Alegre/ RS, 2007.\",\n    \"3)COUPER, JR.; PENNEY, W.R.; FAIR, J.R.; WALAS, S.M. Chemical Process Equipment: Selection and Design. Amsterdam: Elsevier, 2005.\",\n    \"4)TROPEA, C.; YARIN, A.L.; FOSS, J.F. Handbook of Experimental Fluid Mechanics. Ed Springer. Springer-Verlag Berlin Heidelberg. 2007\"\n)\n\n$biblioRange = $d.Paragraphs($biblioIndex).Range\n$biblioRange.Text = [string]::Join([char]11, $references)\n"}
